# B6-PowerPoint.pptx edit - Tue, Jun 16, 2020 11:07:11 AM
#
# 1) Three tables get a different built-in table style applied:
#    {0306776E-D01A-4577-802C-100CA3A21E1E} -> {AD61ED7E-CBBD-40E3-85A1-62D3BD924383}
# 2) The theme actually driving the deck (the one reached through the slide
#    master) swaps palette identity, going from the "Integral" / "Red Violet"
#    colours back to the default "Office Theme" / "Office" colours. We
#    reproduce that by rewriting every slot of the live ThemeColorScheme -
#    PowerPoint serialises those edits straight back into the theme XML part
#    backing the slide master.

$p = $ppt.ActivePresentation

# --- 1) retarget the three tables' style -----------------------------------
$oldStyleId = "{0306776E-D01A-4577-802C-100CA3A21E1E}"
$newStyleId = "{AD61ED7E-CBBD-40E3-85A1-62D3BD924383}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable -and $shape.Table.Style -eq $oldStyleId) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) swap the active theme's colour palette back to the Office defaults -
# RGB values below are COM/OLE colour ints (0xBBGGRR), i.e. the byte-reversed
# form of the target hex colours: 000000,FFFFFF,44546A,E7E6E6,5B9BD5,ED7D31,
# A5A5A5,FFC000,4472C4,70AD47,0563C1,954F72 (dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink)
$officeColors = @{
    1  = 0
    2  = 16777215
    3  = 6968388
    4  = 15132391
    5  = 13998939
    6  = 3243501
    7  = 10855845
    8  = 49407
    9  = 12874308
    10 = 4697456
    11 = 12673797
    12 = 7491477
}

$tcs = $p.Slides.Item(14).ThemeColorScheme
for ($k = 1; $k -le $tcs.Count; $k++) {
    $tcs.Colors($k).RGB = $officeColors[$k]
}
